# New edge weights test sheets
#
# On the "network_optimized_weights" sheet, the gene-name row/column labels
# (ABF1, ACE2, AFT1, CIN5, CUP9, FHL1, GTS1, HAL9, HSF1, MAC1, MSN1, MSN4,
# NRG1, PHD1) are replaced with new synthetic edge labels E-1..E-14 (with
# the first label, ABF1, becoming E-14 and the rest shifting to E-1..E-13).
# This relabels both the header row (B1:O1) and the first data column
# (A2:A15), since the matrix is symmetric-labelled (row i / col i share a
# label). The selection is also updated to span the new header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("network_optimized_weights")
$ws.Activate()

# New labels, in the same left-to-right / top-to-bottom order as the
# existing ABF1..PHD1 labels they replace.
$labels = @("E-14", "E-1", "E-2", "E-3", "E-4", "E-5", "E-6", "E-7", "E-8", "E-9", "E-10", "E-11", "E-12", "E-13")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $label = $labels[$i]

    # Header row, columns B1:O1 (columns 2..15)
    $ws.Cells.Item(1, $i + 2).Value = $label

    # First column, rows A2:A15 (rows 2..15)
    $ws.Cells.Item($i + 2, 1).Value = $label
}

# Matches the updated selection recorded for this sheet.
$ws.Range("B1:O1").Select() | Out-Null
